# Weekly data update: insert a new price record as a new row 267, pushing
# the existing rows 267-337 down to 268-338.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 267 (shifts rows 267:337 -> 268:338).
$ws.Rows.Item(267).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A267").Value = 4
$ws.Range("B267").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C267").Value = "Los Lagos"
$ws.Range("D267").Value = 44508
$ws.Range("E267").Value = 10
$ws.Range("F267").Value = 100112004
$ws.Range("G267").Value = "Cebolla"
$ws.Range("H267").Value = "Sin especificar"
$ws.Range("I267").Value = "Primera"
$ws.Range("J267").Value = 250
$ws.Range("K267").Value = 9000
$ws.Range("L267").Value = 9000
$ws.Range("M267").Value = 9000
$ws.Range("N267").Value = "`$/malla 18 kilos"
$ws.Range("O267").Value = "Perú"
$ws.Range("P267").Value = 500
$ws.Range("Q267").Value = 18
$ws.Range("R267").Value = "Hortaliza"
